$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'230.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'22.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.296"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05599"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.377"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'6.469"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'1.059"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.7823"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1395"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07381"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Value = "'0.02970"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09266"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'0.001673"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'3.251"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04755"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0005791"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.006249"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.005228"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'3.977"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'2.146"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Value = "'0.1051"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'0.0004991"
$ws.Range("D27").Style = "Normal"
$ws.Range("D40").Value = "'0.04028"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.007005"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.003501"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.1039"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.009239"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005441"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.7854"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.04090"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.01010"
$ws.Range("D50").Style = "Normal"
